# Actualización automática 2025-09-09 13:55:09
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("P19").Value = 183.05
$ws1.Range("P34").Value = "1 de 32"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F19").Value = 678.03
$ws2.Range("F34").Value = 3363.05

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D8").Value = 183.05
$ws3.Range("E8").Value = 297.167743214072
$ws3.Range("F8").Value = 0.3811812507694031

$ws3.Range("D15").Value = 3363.05
$ws3.Range("E15").Value = 35379.96881339594
$ws3.Range("F15").Value = 0.08680402568003244
